$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be stored as Text so numeric-looking strings
# (e.g. "0.9995", "1.000") are preserved exactly instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Update Price (D) and Volume(1h) (E) for rows 2-46 ---
$ws.Cells.Item(2, 4).Value = '29.066.65'
$ws.Cells.Item(2, 5).Value = '  -0.42%  '
$ws.Cells.Item(3, 4).Value = '1.827.48'
$ws.Cells.Item(3, 5).Value = '  -0.40%  '
$ws.Cells.Item(4, 4).Value = '0.9995'
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).Value = '240.32'
$ws.Cells.Item(5, 5).Value = '  -0.46%  '
$ws.Cells.Item(6, 4).Value = '0.6144'
$ws.Cells.Item(6, 5).Value = '  -7.71%  '
$ws.Cells.Item(7, 4).Value = '1.001'
$ws.Cells.Item(7, 5).Value = '  +0.12%  '
$ws.Cells.Item(8, 4).Value = '44.47'
$ws.Cells.Item(8, 5).Value = '  +6.10%  '
$ws.Cells.Item(9, 4).Value = '0.07308'
$ws.Cells.Item(9, 5).Value = '  -1.37%  '
$ws.Cells.Item(10, 4).Value = '0.2910'
$ws.Cells.Item(10, 5).Value = '  -0.84%  '
$ws.Cells.Item(11, 4).Value = '22.56'
$ws.Cells.Item(11, 5).Value = '  -0.33%  '
$ws.Cells.Item(12, 4).Value = '0.07693'
$ws.Cells.Item(12, 5).Value = '  -0.39%  '
$ws.Cells.Item(13, 4).Value = '1.826.54'
$ws.Cells.Item(13, 5).Value = '  -3.00%  '
$ws.Cells.Item(14, 4).Value = '4.956'
$ws.Cells.Item(14, 5).Value = '  -0.69%  '
$ws.Cells.Item(15, 4).Value = '0.6586'
$ws.Cells.Item(15, 5).Value = '  -1.59%  '
$ws.Cells.Item(16, 4).Value = '81.66'
$ws.Cells.Item(16, 5).Value = '  -1.56%  '
$ws.Cells.Item(17, 4).Value = '0.000008913'
$ws.Cells.Item(17, 5).Value = '  +5.94%  '
$ws.Cells.Item(18, 4).Value = '6.017'
$ws.Cells.Item(18, 5).Value = '  -1.58%  '
$ws.Cells.Item(19, 4).Value = '29.059.19'
$ws.Cells.Item(19, 5).Value = '  -0.86%  '
$ws.Cells.Item(20, 4).Value = '2.076.12'
$ws.Cells.Item(20, 5).Value = '  -3.59%  '
$ws.Cells.Item(21, 4).Value = '224.63'
$ws.Cells.Item(21, 5).Value = '  -0.90%  '
$ws.Cells.Item(22, 4).Value = '12.35'
$ws.Cells.Item(22, 5).Value = '  -1.03%  '
$ws.Cells.Item(23, 4).Value = '1.001'
$ws.Cells.Item(23, 5).Value = '  -0.01%  '
$ws.Cells.Item(24, 4).Value = '7.116'
$ws.Cells.Item(24, 5).Value = '  -1.11%  '
$ws.Cells.Item(25, 4).Value = '1.001'
$ws.Cells.Item(25, 5).Value = '  +0.06%  '
$ws.Cells.Item(26, 4).Value = '159.01'
$ws.Cells.Item(26, 5).Value = '  -0.41%  '
$ws.Cells.Item(27, 4).Value = '8.417'
$ws.Cells.Item(27, 5).Value = '  -2.40%  '
$ws.Cells.Item(28, 4).Value = '0.1339'
$ws.Cells.Item(28, 5).Value = '  -4.98%  '
$ws.Cells.Item(29, 4).Value = '17.73'
$ws.Cells.Item(29, 5).Value = '  -1.31%  '
$ws.Cells.Item(30, 4).Value = '1.497'
$ws.Cells.Item(30, 5).Value = '  -1.06%  '
$ws.Cells.Item(31, 4).Value = '4.052'
$ws.Cells.Item(31, 5).Value = '  -1.39%  '
$ws.Cells.Item(32, 4).Value = '4.024'
$ws.Cells.Item(32, 5).Value = '  -0.67%  '
$ws.Cells.Item(33, 4).Value = '1.197'
$ws.Cells.Item(33, 5).Value = '  +1.24%  '
$ws.Cells.Item(34, 4).Value = '0.05277'
$ws.Cells.Item(34, 5).Value = '  -0.41%  '
$ws.Cells.Item(35, 4).Value = '1.834'
$ws.Cells.Item(35, 5).Value = '  -2.21%  '
$ws.Cells.Item(36, 4).Value = '1.144'
$ws.Cells.Item(36, 5).Value = '  +0.88%  '
$ws.Cells.Item(37, 4).Value = '0.7290'
$ws.Cells.Item(37, 5).Value = '  -4.25%  '
$ws.Cells.Item(38, 4).Value = '2.652'
$ws.Cells.Item(38, 5).Value = '  -0.77%  '
$ws.Cells.Item(39, 4).Value = '1.284.56'
$ws.Cells.Item(39, 5).Value = '  +1.01%  '
$ws.Cells.Item(40, 4).Value = '0.01782'
$ws.Cells.Item(40, 5).Value = '  -0.79%  '
$ws.Cells.Item(41, 4).Value = '2.742'
$ws.Cells.Item(41, 5).Value = '  +0.68%  '
$ws.Cells.Item(42, 4).Value = '6.354'
$ws.Cells.Item(42, 5).Value = '  +6.38%  '
$ws.Cells.Item(43, 4).Value = '0.9003'
$ws.Cells.Item(43, 5).Value = '  -2.98%  '
$ws.Cells.Item(44, 4).Value = '1.000'
$ws.Cells.Item(44, 5).Value = '  -0.09%  '
$ws.Cells.Item(45, 4).Value = '101.83'
$ws.Cells.Item(45, 5).Value = '  -1.01%  '
$ws.Cells.Item(46, 4).Value = '1.975.23'
$ws.Cells.Item(46, 5).Value = '  -3.50%  '

# --- Rows 47-51 reshuffled: update Coin, Link, Price, Volume(1h) ---
$ws.Cells.Item(47, 2).Value = 'Mantle'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(47, 4).Value = '0.5125'
$ws.Cells.Item(47, 5).Value = '  -0.70%  '
$ws.Cells.Item(48, 2).Value = 'Aave'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(48, 4).Value = '63.81'
$ws.Cells.Item(48, 5).Value = '  +0.60%  '
$ws.Cells.Item(49, 2).Value = 'RenderToken'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(49, 4).Value = '1.709'
$ws.Cells.Item(49, 5).Value = '  -3.79%  '
$ws.Cells.Item(50, 2).Value = 'TheSandbox'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(50, 4).Value = '0.3956'
$ws.Cells.Item(50, 5).Value = '  -2.15%  '
$ws.Cells.Item(51, 2).Value = 'Cronos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(51, 4).Value = '0.05791'
$ws.Cells.Item(51, 5).Value = '  -2.32%  '
